# Ratio constraints implemented and tested.
# Updates the "factors" ratio values, refreshes the RSE threshold column (D),
# clears the stray D9 value, and moves the active-sheet/selection state from
# RSE -> factors (B2:B7 selected) while RSE's own selection settles on G5.

$wb = $excel.ActiveWorkbook

$factors = $wb.Worksheets.Item("factors")
$rse = $wb.Worksheets.Item("RSE")

# --- factors sheet: updated ratio values (column B) ---
$factors.Range("B2").Value = 0.9
$factors.Range("B3").Value = 0.3
$factors.Range("B4").Value = 0.1
$factors.Range("B5").Value = 0.5
$factors.Range("B6").Value = 0.1
$factors.Range("B7").Value = 1.5

# --- RSE sheet: column D updates ---
$rse.Range("D2").Value = 1000000
$rse.Range("D3").Value = 1000000
$null = $rse.Range("D9").ClearContents()

# --- view / selection state ---
# RSE no longer the active tab; its own selection moves to G5.
$null = $rse.Range("G5").Select()

# factors becomes the active tab, with B2:B7 selected.
$null = $factors.Activate()
$null = $factors.Range("B2:B7").Select()
